# Generate Report for Handoff
# Updates the "Ready for handoff" rows (4-7) on the zh-cn and de-de sheets:
#   - Priority column (E) changes from "low" to "ht"
#   - Latest Handoff Datetime column (H) is refreshed to a new generation time

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# zh-cn: rows 4-7 get new handoff datetime 2016-09-04 08:34:27
$zhcn.Range("E4:E7").Value = "ht"
$zhcn.Range("H4:H7").Value = "2016-09-04 08:34:27"

# de-de: rows 4-7 get new handoff datetime 2016-09-04 08:34:32
$dede.Range("E4:E7").Value = "ht"
$dede.Range("H4:H7").Value = "2016-09-04 08:34:32"

# Overview sheet: "Latest HO Xliff Generate Date" for the "Ready for handoff" rows (4-7)
# shares the same updated timestamp as the de-de handoff time
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G4:G7").Value = "2016-09-04 08:34:32"
